$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'questions = [
    {
        "title": "You work for a well-established ecommerce business that has gathered a lot of data from existing customers. You are planning a new reach campaign in which you intend to target both existing customers and newer audiences, with a 50-50 budget split.  What campaign structure should you use?",
        "ques_type": 2,
        "options": [
            "Two separate campaigns, with the total budget set at the account level.",
            "One campaign with one ad set, using the total budget to target both audiences.",
            "One campaign with two separate ad sets, with the budget set at the campaign level.",
            "One campaign with two separate ad sets, with the budget set at the ad set level."
        ],
        "score": "One campaign with two separate ad sets, with the budget set at the ad set level."
    },
    {
        "title": "You are the digital ads manager at a beverage company. You are running an engagement campaign that is optimized for video views. The campaign is working fairly well, but you now want to update it to narrow its targeting.How should you do this?",
        "ques_type": 2,
        "options": [
            "Edit the audience &gt Narrow audience &gt Update",
            "Edit the audience &gt Define further &gt Publish",
            "Duplicate the campaign &gt Narrow audience &gt Save as new",
            "Duplicate the ad &gt Narrow audience &gt Update"
        ],
        "score": "Edit the audience &gt Define further &gt Publish"
    },
    {
        "title": "You work for a makeup brand. You want to run a video carousel on Advantage+ placements. Which creative assets should you request?",
        "ques_type": 2,
        "options": [
            "15 videos with a ratio of 9:16",
            "15 videos with a ratio of 1:1",
            "10 videos with a ratio of 1:1",
            "10 videos with a ratio of 16:9"
        ],
        "score": "10 videos with a ratio of 1:1"
    },
    {
        "title": "You work for a performance marketing agency and are evaluating past campaigns.  Which of the following campaigns has the highest return on ad spend (ROAS = revenue/spend)?",
        "ques_type": 2,
        "options": [
            "Campaign #1",
            "Campaign #2",
            "Campaign #3",
            "Campaign #4"
        ],
        "score": "Campaign #2"
    }
]'

# The sheet originally has two rows: A1 (numeric 0, bold+bordered style) and
# A2 (shared string with the questions blob). The edit collapses this down
# to a single A1 cell holding the (now pretty-printed) questions text with
# plain/default formatting.
$ws.Range("A2").EntireRow.Delete()
$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"
$ws.Rows.Item(1).AutoFit()
